$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.909.59"
$ws.Range("E2").Value = "  -3.65%  "
$ws.Range("D3").Value = "'3.504.29"
$ws.Range("E3").Value = "  -2.86%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'586.61"
$ws.Range("E5").Value = "  -3.25%  "
$ws.Range("D6").Value = "'132.24"
$ws.Range("E6").Value = "  -5.51%  "
$ws.Range("D7").Value = "'3.503.77"
$ws.Range("E7").Value = "  -2.83%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.494"
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").Value = "'0.124"
$ws.Range("E10").Value = "  -2.50%  "
$ws.Range("D11").Value = "'7.13"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").Value = "'0.385"
$ws.Range("E12").Value = "  -2.47%  "
$ws.Range("D13").Value = "'4.101.99"
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("D14").Value = "'27.87"
$ws.Range("E14").Value = "  -2.62%  "
$ws.Range("E15").Value = "  -4.58%  "
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "'3.502.60"
$ws.Range("E17").Value = "  -2.90%  "
$ws.Range("D18").Value = "'64.044.71"
$ws.Range("E18").Value = "  -3.53%  "
$ws.Range("D19").Value = "'10.04"
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("D20").Value = "'14.48"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("D21").Value = "'5.69"
$ws.Range("E21").Value = "  -4.15%  "
$ws.Range("D22").Value = "'391.59"
$ws.Range("E22").Value = "  -2.13%  "
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("D24").Value = "'3.645.93"
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("D25").Value = "'72.95"
$ws.Range("E25").Value = "  -3.09%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -6.96%  "
$ws.Range("E28").Value = "  -3.55%  "
$ws.Range("D29").Value = "'7.50"
$ws.Range("E29").Value = "  -8.54%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("E32").Value = "  -5.04%  "
$ws.Range("D33").Value = "'3.512.23"
$ws.Range("E33").Value = "  -2.86%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'23.89"
$ws.Range("E35").Value = "  -3.42%  "
$ws.Range("D36").Value = "'0.144"
$ws.Range("E36").Value = "  -4.16%  "
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("E38").Value = "  -4.41%  "
$ws.Range("D39").Value = "'6.95"
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("D40").Value = "'166.99"
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("E41").Value = "  -4.44%  "
$ws.Range("D42").Value = "'27.27"
$ws.Range("E42").Value = "  +3.30%  "
$ws.Range("E43").Value = "  -3.83%  "
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "'41.92"
$ws.Range("E45").Value = "  -3.04%  "
$ws.Range("D46").Value = "'1.20"
$ws.Range("E46").Value = "  -6.03%  "
$ws.Range("E47").Value = "  -4.29%  "
$ws.Range("E48").Value = "  -5.39%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'6.90"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "'2.447.81"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").Value = "'0.902"
$ws.Range("E51").Value = "  -1.37%  "
